$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: mark as duplicate-fixed TRUE with corrected count
$ws.Range("A2").Value = $true
$ws.Range("B2").Value = 9988

# Remove the now-obsolete row 3 (previously held TRUE/9994, folded into row 2)
$ws.Rows("3").Delete()
